# Generate Report for Handback
# Updates the handoff/handback timestamps recorded for the
# 28cb43ac... and 45dc5af5... entries on the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-20 00:16:59"
$zhcn.Range("H2").Value = "2016-03-20 00:17:20"
$zhcn.Range("E3").Value = "2016-03-20 00:16:59"
$zhcn.Range("H3").Value = "2016-03-20 00:17:20"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-20 00:17:02"
$dede.Range("H2").Value = "2016-03-20 00:17:25"
$dede.Range("E3").Value = "2016-03-20 00:17:02"
$dede.Range("H3").Value = "2016-03-20 00:17:25"
